$wb = $excel.ActiveWorkbook

# --- Misc sheet: "version"/0.2 -> "model"/"Linear" with a dropdown list ---
$ws = $wb.Worksheets.Item("Misc")
$ws.Range("A1").Value = "model"
$ws.Range("B1").Value = "Linear"
$ws.Range("B1").Validation.Add(3, 1, 1, """Linear, 2FI, Quadratic""") | Out-Null

# Make "Misc" the active sheet/tab with its own selection, matching the
# restored state (tab 4 selected, cell I6 highlighted).
$ws.Activate() | Out-Null
$ws.Range("I6").Select() | Out-Null
